$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at row 51 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 51-119 down to 52-120 (their values are preserved
# verbatim by the native row insert, which matches the diff: each old row N's
# data reappears unchanged as row N+1, and the former last row (119) becomes
# the new last row (120)).
$ws.Rows(51).Insert()

# Populate the newly inserted row 51 with this week's record.
$ws.Range("A51").Value = 4
$ws.Range("B51").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C51").Value = "Los Lagos"
$ws.Range("D51").Value = 45118
$ws.Range("E51").Value = 10
$ws.Range("F51").Value = 100112043
$ws.Range("G51").Value = "Pepino dulce"
$ws.Range("H51").Value = "Cultivar IV Región"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 60
$ws.Range("K51").Value = 21000
$ws.Range("L51").Value = 21000
$ws.Range("M51").Value = 21000
$ws.Range("N51").Value = "$/bandeja 18 kilos"
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 1167
$ws.Range("Q51").Value = 18
$ws.Range("R51").Value = "Hortaliza"
